$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers for columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header formatting (bold, centered, bordered) used by the
# existing header row (copy format from H1, which already has it).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Fill the data rows: I is always 1, J mirrors the value already in H.
for ($r = 2; $r -le 35; $r++) {
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $ws.Cells.Item($r, 8).Value2
}

$excel.CutCopyMode = 0
